$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 4969.3335
$ws.Range("J48").Value = 4969.3335
$ws.Range("L48").Value = 14908.0005
$ws.Range("N48").Value = -15492.0005

$ws.Range("H56").Value = 4969.3335
$ws.Range("J56").Value = 4969.3335
$ws.Range("L56").Value = 14908.0005
$ws.Range("N56").Value = -15976.0005

$ws.Range("H99").Value = 3103.0908
$ws.Range("J99").Value = 10166
$ws.Range("L99").Value = 30498
$ws.Range("N99").Value = -33494

$ws.Range("H107").Value = 1372.5238
$ws.Range("J107").Value = 971.6667
$ws.Range("L107").Value = 971.6667
$ws.Range("N107").Value = -4811.6667

$ws.Range("H137").Value = 6444.614
$ws.Range("I137").Value = 3880.6858
$ws.Range("K137").Value = 11642.0574
$ws.Range("M137").Value = -9092.057400000002

$ws.Range("H138").Value = 4311.5
$ws.Range("J138").Value = 5678.1934
$ws.Range("L138").Value = 17034.5802
$ws.Range("N138").Value = -27314.5802

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3325
$ws.Range("I45").Value = 3325
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 3325
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -2948
$ws.Range("N45").ClearContents()

$ws.Range("H63").Value = 300
$ws.Range("I63").Value = 275
$ws.Range("J63").Value = 325
$ws.Range("K63").Value = 275
$ws.Range("L63").Value = 325
$ws.Range("M63").Value = 411
$ws.Range("N63").Value = -1697

$ws.Range("H66").Value = 300
$ws.Range("I66").Value = 275
$ws.Range("J66").Value = 325
$ws.Range("K66").Value = 1375
$ws.Range("L66").Value = 1625
$ws.Range("M66").Value = 2057
$ws.Range("N66").Value = -8489

$ws.Range("H74").Value = 4018.6365
$ws.Range("I74").Value = 4569.625
$ws.Range("K74").Value = 4569.625
$ws.Range("M74").Value = -3695.625

$ws.Range("H77").Value = 4018.6365
$ws.Range("I77").Value = 4569.625
$ws.Range("K77").Value = 22848.125
$ws.Range("M77").Value = -18480.125

$ws.Range("H102").Value = 8707.357
$ws.Range("I102").Value = 8817.333000000001
$ws.Range("K102").Value = 8817.333000000001
$ws.Range("M102").Value = -7195.333000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3043.9
$ws.Range("I105").Value = 2521
$ws.Range("K105").Value = 2521
$ws.Range("M105").Value = -774

$ws.Range("H134").Value = 1732736.4
$ws.Range("I134").Value = 2090958.5
$ws.Range("J134").Value = 13270
$ws.Range("K134").Value = 6272875.5
$ws.Range("L134").Value = 39810
$ws.Range("M134").Value = -6270340.5
$ws.Range("N134").Value = -44880

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5641.346
$ws.Range("I31").Value = 2272.6
$ws.Range("K31").Value = 2272.6
$ws.Range("M31").Value = -1977.6

$ws.Range("H34").Value = 5641.346
$ws.Range("I34").Value = 2272.6
$ws.Range("K34").Value = 2272.6
$ws.Range("M34").Value = -2070.6

$ws.Range("H58").Value = 28577938
$ws.Range("I58").Value = 47624020
$ws.Range("J58").Value = 8812.214
$ws.Range("K58").Value = 47624020
$ws.Range("L58").Value = 8812.214
$ws.Range("M58").Value = -47623817
$ws.Range("N58").Value = -9218.214

$ws.Range("H70").Value = 42000
$ws.Range("J70").Value = 45000
$ws.Range("L70").Value = 45000
$ws.Range("N70").Value = -45630

$ws.Range("H73").Value = 42000
$ws.Range("J73").Value = 45000
$ws.Range("L73").Value = 45000
$ws.Range("N73").Value = -47184

$ws.Range("H99").Value = 3835350
$ws.Range("I99").Value = 11114143
$ws.Range("J99").Value = 4406.3687
$ws.Range("K99").Value = 11114143
$ws.Range("L99").Value = 4406.3687
$ws.Range("M99").Value = -11112645
$ws.Range("N99").Value = -7402.3687

$ws.Range("H111").Value = 74331
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 74331
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 74331
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -82511

$ws.Range("H126").Value = 3835350
$ws.Range("I126").Value = 11114143
$ws.Range("J126").Value = 4406.3687
$ws.Range("K126").Value = 33342429
$ws.Range("L126").Value = 13219.1061
$ws.Range("M126").Value = -33339959
$ws.Range("N126").Value = -18159.1061

$ws.Range("H134").Value = 30307808
$ws.Range("I134").Value = 41670676
$ws.Range("K134").Value = 125012028
$ws.Range("M134").Value = -125009493

$ws.Range("H136").Value = 28577938
$ws.Range("I136").Value = 47624020
$ws.Range("J136").Value = 8812.214
$ws.Range("K136").Value = 142872060
$ws.Range("L136").Value = 26436.642
$ws.Range("M136").Value = -142869510
$ws.Range("N136").Value = -31536.642

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 626221.4
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 626221.4
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1878664.2
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -1883004.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 30026.375
$ws.Range("J39").Value = 32173
$ws.Range("L39").Value = 32173
$ws.Range("N39").Value = -33237

$ws.Range("H102").Value = 4467.259
$ws.Range("I102").Value = 3254
$ws.Range("J102").Value = 6452.591
$ws.Range("K102").Value = 3254
$ws.Range("L102").Value = 6452.591
$ws.Range("M102").Value = -1632
$ws.Range("N102").Value = -9696.591

$ws.Range("H126").Value = 20005132
$ws.Range("I126").Value = 31253160
$ws.Range("J126").Value = 8634.777
$ws.Range("K126").Value = 93759480
$ws.Range("L126").Value = 25904.331
$ws.Range("M126").Value = -93757010
$ws.Range("N126").Value = -30844.331

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 36669.5

$ws.Range("H20").Value = 11461.462
$ws.Range("J20").Value = 11461.462
$ws.Range("L20").Value = 11461.462
$ws.Range("N20").Value = -11913.462

$ws.Range("H22").Value = 1040.2858
$ws.Range("J22").Value = 939.4286
$ws.Range("L22").Value = 939.4286
$ws.Range("N22").Value = -1529.4286

$ws.Range("H27").Value = 1040.2858
$ws.Range("J27").Value = 939.4286
$ws.Range("L27").Value = 939.4286
$ws.Range("N27").Value = -1153.4286

$ws.Range("H40").Value = 4332.3076
$ws.Range("I40").Value = 3610.842
$ws.Range("J40").Value = 6290.5713
$ws.Range("K40").Value = 3610.842
$ws.Range("L40").Value = 6290.5713
$ws.Range("M40").Value = -3474.842
$ws.Range("N40").Value = -6562.5713

$ws.Range("H46").Value = 31251024
$ws.Range("I46").Value = 1106.6666
$ws.Range("J46").Value = 50000976
$ws.Range("K46").Value = 1106.6666
$ws.Range("L46").Value = 50000976
$ws.Range("M46").Value = -918.6666
$ws.Range("N46").Value = -50001352

$ws.Range("H68").Value = 1841.9474
$ws.Range("I68").Value = 1466.4667
$ws.Range("J68").Value = 3250
$ws.Range("K68").Value = 1466.4667
$ws.Range("L68").Value = 3250
$ws.Range("M68").Value = -717.4666999999999
$ws.Range("N68").Value = -4748

$ws.Range("H71").Value = 1841.9474
$ws.Range("I71").Value = 1466.4667
$ws.Range("J71").Value = 3250
$ws.Range("K71").Value = 7332.3335
$ws.Range("L71").Value = 16250
$ws.Range("M71").Value = -3588.3335
$ws.Range("N71").Value = -23738

$ws.Range("H93").Value = 1531.8334
$ws.Range("I93").Value = 1593.3334
$ws.Range("J93").Value = 1408.8334
$ws.Range("K93").Value = 1593.3334
$ws.Range("L93").Value = 1408.8334
$ws.Range("M93").Value = -345.3334
$ws.Range("N93").Value = -3904.8334

$ws.Range("H122").Value = 12154.077
$ws.Range("I122").Value = 21667.166
$ws.Range("K122").Value = 65001.49800000001
$ws.Range("M122").Value = -62551.49800000001

$ws.Range("H132").Value = 4846.289
$ws.Range("I132").Value = 3942
$ws.Range("K132").Value = 11826
$ws.Range("M132").Value = -9296

$ws.Range("H136").Value = 27782882
$ws.Range("I136").Value = 33338364
$ws.Range("K136").Value = 100015092
$ws.Range("M136").Value = -100012542

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 18403.857
$ws.Range("I62").Value = 14019.5
$ws.Range("J62").Value = 24249.666
$ws.Range("K62").Value = 14019.5
$ws.Range("L62").Value = 24249.666
$ws.Range("M62").Value = -13395.5
$ws.Range("N62").Value = -25497.666

$ws.Range("H65").Value = 18403.857
$ws.Range("I65").Value = 14019.5
$ws.Range("J65").Value = 24249.666
$ws.Range("K65").Value = 70097.5
$ws.Range("L65").Value = 121248.33
$ws.Range("M65").Value = -66977.5
$ws.Range("N65").Value = -127488.33

$ws.Range("H107").Value = 4836.143
$ws.Range("I107").Value = 757.7692
$ws.Range("K107").Value = 2273.3076
$ws.Range("M107").Value = -353.3076000000001

$ws.Range("H122").Value = 2154.6667
$ws.Range("I122").Value = 1849.3636
$ws.Range("J122").Value = 2994.25
$ws.Range("K122").Value = 5548.0908
$ws.Range("L122").Value = 8982.75
$ws.Range("M122").Value = -3098.0908
$ws.Range("N122").Value = -13882.75

$ws.Range("H132").Value = 5587.1143
$ws.Range("I132").Value = 5097.636
$ws.Range("J132").Value = 6415.4614
$ws.Range("K132").Value = 15292.908
$ws.Range("L132").Value = 19246.3842
$ws.Range("M132").Value = -12762.908
$ws.Range("N132").Value = -24306.3842
